$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-13 Saturday" "2025-09-14 Sunday"
Replace-Text "246÷2=" "509÷7="
Replace-Text "908÷3=" "101÷8="
Replace-Text "881÷6=" "756÷8="
Replace-Text "506÷6=" "504÷4="
Replace-Text "843÷7=" "111÷9="
Replace-Text "909÷9=" "601÷7="
Replace-Text "198÷8=" "389÷2="
Replace-Text "579÷9=" "997÷3="
Replace-Text "782÷5=" "335÷8="
Replace-Text "105÷4=" "198÷5="
Replace-Text "762÷3=" "985÷8="
Replace-Text "448÷9=" "185÷6="
Replace-Text "216÷5=" "751÷6="
Replace-Text "167÷7=" "350÷6="
Replace-Text "592÷4=" "455÷6="
Replace-Text "331÷7=" "609÷8="
Replace-Text "675÷2=" "365÷3="
Replace-Text "436÷6=" "427÷6="
Replace-Text "946÷9=" "871÷3="
Replace-Text "477÷6=" "173÷7="
Replace-Text "119÷6=" "292÷3="
Replace-Text "185÷3=" "502÷3="
Replace-Text "578÷9=" "102÷4="
Replace-Text "433÷3=" "732÷6="
Replace-Text "653÷6=" "392÷6="
